$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the old "QC" column (current column D), shifting
# "QC"/"QT" one column to the right, and give the new column the same
# formatting (width/style) as the "Define" column (C).
$ws.AutoFilterMode = $false
$ws.Columns("D").Insert()
$ws.Columns("D").ColumnWidth = 50.7109375
$ws.Range("D1").Value = "Note"

# Re-apply the AutoFilter over the new, wider header range.
$ws.Range("A1:F1").AutoFilter(1, "x")
$ws.Range("A1:F1").AutoFilter(1)

# Keep the hidden _FilterDatabase defined name in sync with the AutoFilter range.
$wb.Names.Item(1).RefersTo = "='iciba-collins'!`$A`$1:`$F`$1"
